$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 77: change to "Nacional vs Santa Clara" data (was "Benfica B vs Penafiel") ---
$ws.Range("F77").Value = "Nacional"
$ws.Range("H77").Value = "Santa Clara"
$ws.Range("I77").Value = 1
$ws.Range("J77").Value = 2.98
$ws.Range("K77").Value = "01/11/2023 16:12"
$ws.Range("L77").Value = 2.81
$ws.Range("M77").Value = "04/11/2023 18:58"
$ws.Range("N77").Value = 3.27
$ws.Range("O77").Value = "01/11/2023 16:12"
$ws.Range("P77").Value = 3.23
$ws.Range("Q77").Value = "04/11/2023 18:52"
$ws.Range("R77").Value = 2.39
$ws.Range("S77").Value = "01/11/2023 16:12"
$ws.Range("T77").Value = 2.72
$ws.Range("U77").Value = "04/11/2023 18:52"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/nacional-santa-clara/xQH2R8bH/"

# --- Row 78: change to "Benfica B vs Penafiel" data (was "Nacional vs Santa Clara") ---
$ws.Range("F78").Value = "Benfica B"
$ws.Range("H78").Value = "Penafiel"
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1.91
$ws.Range("K78").Value = "29/10/2023 16:42"
$ws.Range("L78").Value = 2.26
$ws.Range("M78").Value = "04/11/2023 18:53"
$ws.Range("N78").Value = 3.67
$ws.Range("O78").Value = "29/10/2023 16:42"
$ws.Range("P78").Value = 3.59
$ws.Range("Q78").Value = "04/11/2023 18:53"
$ws.Range("R78").Value = 4.01
$ws.Range("S78").Value = "29/10/2023 16:42"
$ws.Range("T78").Value = 3.2
$ws.Range("U78").Value = "04/11/2023 18:53"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/benfica-penafiel/xjmbUAEb/"

# --- New row 83: Feirense vs Os Belenenses ---
# Copy formats from row 82 (last existing data row) so the new row matches
# the workbook's styling (bold/border on column A, date format on column E).
$ws.Range("A82:V82").Copy()
$ws.Range("A83:V83").PasteSpecial(-4122)

$ws.Range("A83").Value = 82
$ws.Range("B83").Value = "portugal"
$ws.Range("C83").Value = "liga-portugal-2"
$ws.Range("D83").Value = "2023-2024"
$ws.Range("E83").Value = 45240.79166666666
$ws.Range("F83").Value = "Feirense"
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = "Os Belenenses"
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2.33
$ws.Range("K83").Value = "07/11/2023 22:17"
$ws.Range("L83").Value = 2.02
$ws.Range("M83").Value = "10/11/2023 18:23"
$ws.Range("N83").Value = 3.41
$ws.Range("O83").Value = "07/11/2023 22:17"
$ws.Range("P83").Value = 3.44
$ws.Range("Q83").Value = "10/11/2023 18:23"
$ws.Range("R83").Value = 2.95
$ws.Range("S83").Value = "07/11/2023 22:17"
$ws.Range("T83").Value = 4.02
$ws.Range("U83").Value = "10/11/2023 18:23"
$ws.Range("V83").Value = "https://www.betexplorer.com/football/portugal/liga-portugal-2/feirense-cf-os-belenenses/OGjLqjEo/"
